# Edit script: add the 'Tienda' (HUJZ-001..003) user-story paragraphs
# after the existing 'HULJ-003 Registro de videojuegos' story, replacing
# two of the three trailing blank paragraphs with the new content.

$d = $word.ActiveDocument

# Locate the anchor paragraph (end of the HULJ-003 story) by scanning for
# its unique trailing text, then walk forward with Paragraph.Next() so we
# do not depend on absolute paragraph indices (robust to minor shifts).
$paras = $d.Paragraphs
$anchorIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like "*evitar conflictos de intereses y regal*") {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "Anchor paragraph (HULJ-003 story) not found"
}
$anchorPara = $paras.Item($anchorIndex)

$blank1 = $anchorPara.Next()
$blank2 = $blank1.Next()
$blank3 = $blank2.Next()

if (($blank1.Range.Text -ne [string][char]13) -or ($blank2.Range.Text -ne [string][char]13) -or ($blank3.Range.Text -ne [string][char]13)) {
    throw "Expected three blank paragraphs after the anchor paragraph"
}

# Delete the 3rd blank paragraph entirely (merges it away), leaving the
# 2nd blank paragraph as the single insertion point for the new content.
$blank3.Range.Delete() | Out-Null

# Re-fetch the remaining blank paragraph; indices/ranges shift after delete.
$blank1b = $anchorPara.Next()
$target = $blank1b.Next()
$targetRange = $target.Range

# Build the WordprocessingML package fragment containing the six new
# paragraphs (the HUJZ-001/002/003 stories) plus one trailing blank
# paragraph, exactly mirroring the target OOXML.
$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$xmlBody = @'
    <w:p>
      <w:pPr>
        <w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:noProof/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>HUJZ-001 Administracion de sistema</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:noProof/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>Como</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> administrador del sistema,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>quiero</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> poder registrar nuevas tiendas con su información básica (nombre, dirección, contacto),</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>para</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> gestionar múltiples sucursales dentro del sistema de ventas.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:noProof/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>HUJZ-002 Administracion de sistema</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>Visualización por tienda</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>Quiero:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> gerente regional,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>quiero</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> ver un reporte de ventas filtrado por tienda,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>para</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> analizar el rendimiento individual de cada sucursal y tomar decisiones estratégicas.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:noProof/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:noProof/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>HUJZ</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:noProof/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>-003 Asignaci</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:noProof/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t>on a productor por</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:bCs/>
          <w:noProof/>
          <w:color w:val="4472C4" w:themeColor="accent1"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> tienda </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>Como</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> encargado de inventario,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>quiero</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> asignar productos específicos a cada tienda,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:bCs/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>para</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>asegurar que el stock disponible esté correctamente distribuido entre las sucursales.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'@
$xmlFragment = $xmlHeader + $xmlBody + $xmlFooter

$targetRange.InsertXML($xmlFragment) | Out-Null

Write-Host "Paragraphs after edit: $($d.Paragraphs.Count)"
